$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 778, shifting rows 778:844 down to 779:845.
$ws.Rows(778).Insert()

# Populate the newly inserted row 778 with the new record's data.
$ws.Cells.Item(778, 1).Value = 10
$ws.Cells.Item(778, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(778, 3).Value = "La Araucanía"
$ws.Cells.Item(778, 4).Value = 45106
$ws.Cells.Item(778, 5).Value = 9
$ws.Cells.Item(778, 6).Value = 100112032
$ws.Cells.Item(778, 7).Value = "Zapallo italiano"
$ws.Cells.Item(778, 8).Value = "Bola 8"
$ws.Cells.Item(778, 9).Value = "Primera"
$ws.Cells.Item(778, 10).Value = 55
$ws.Cells.Item(778, 11).Value = 16000
$ws.Cells.Item(778, 12).Value = 16000
$ws.Cells.Item(778, 13).Value = 16000
$ws.Cells.Item(778, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(778, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(778, 16).Value = 320
$ws.Cells.Item(778, 17).Value = 50
$ws.Cells.Item(778, 18).Value = "Hortaliza"
